# Regenerate merged AHB files
# - rename the "_old" / "_new" header suffixes to "_FV2310" / "_FV2404"
# - wrap the data range in a table (Table1)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $colOld = $i + 1       # columns A..J
    $colNew = $i + 12      # columns L..U (K = "diff", unchanged)
    $ws.Cells.Item(1, $colOld).Value = $baseNames[$i] + "_FV2310"
    $ws.Cells.Item(1, $colNew).Value = $baseNames[$i] + "_FV2404"
}

# Turn the used range into a table
$rng = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
